# Edit: "Added extra tests for speed saturation"
# Adds a TEST 6b variant table next to TEST 6a, tweaks several test result
# values, and appends four new test sections (TEST 12 x2, TEST 13 x2).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Simple numeric / value corrections in existing rows
# ---------------------------------------------------------------------
$ws.Range("N37").Value = 4
$ws.Range("N45").Value = 3
$ws.Range("M52").Value = 0
$ws.Range("M59").Value = 0
$ws.Range("G60").Value = 20
$ws.Range("I60").Value = 35
$ws.Range("L60").Value = 35
$ws.Range("M60").Value = 2
$ws.Range("N66").Value = 4
$ws.Range("M79").Value = 0
$ws.Range("I80").Value = 34
$ws.Range("L80").Value = 34
$ws.Range("L87").Value = 37.5
$ws.Range("M87").Value = 20.28
$ws.Range("H94").Value = 0
$ws.Range("L94").Value = 32.5
$ws.Range("M94").Value = 0

# ---------------------------------------------------------------------
# 2. Rename TEST 6 header -> TEST 6a, add a parallel TEST 6b table
# ---------------------------------------------------------------------
$ws.Range("A55").Value = "TEST 6a - Moving from ON to STBY to DISABLE with resume button"

$hdr = $ws.Range("P55:X55")
$hdr.Merge()
$hdr.Font.Bold = $true
$hdr.HorizontalAlignment = -4108
$ws.Range("P55").Value = "TEST 6b - Moving from ON to STBY to DISABLE with resume button"

$ws.Range("Y55").Font.Bold = $true
$ws.Range("Z55").Font.Bold = $true

$ws.Range("AA55").Value = "CruiseSpeed "
$ws.Range("AB55").Value = "ThrottleCmd "
$ws.Range("AC55").Value = "CruiseState"

$ws.Range("P57").Value = "On"
$ws.Range("Q57").Value = "Off "
$ws.Range("R57").Value = "Resume "
$ws.Range("S57").Value = "Set"
$ws.Range("T57").Value = "QuickAccel"
$ws.Range("U57").Value = "QuickDecel"
$ws.Range("V57").Value = " Accel"
$ws.Range("W57").Value = "Break"
$ws.Range("X57").Value = "Speed"

$ws.Range("P58").Value = $true
$ws.Range("Q58").Value = $false
$ws.Range("R58").Value = $false
$ws.Range("S58").Value = $false
$ws.Range("T58").Value = $false
$ws.Range("U58").Value = $false
$ws.Range("V58").Value = 0
$ws.Range("W58").Value = 0
$ws.Range("X58").Value = 35
$ws.Range("AA58").Value = 35
$ws.Range("AB58").Value = 0
$ws.Range("AC58").Value = 2

$ws.Range("P59").Value = $false
$ws.Range("Q59").Value = $false
$ws.Range("R59").Value = $false
$ws.Range("S59").Value = $false
$ws.Range("T59").Value = $false
$ws.Range("U59").Value = $false
$ws.Range("V59").Value = 0
$ws.Range("W59").Value = 20
$ws.Range("X59").Value = 35
$ws.Range("AA59").Value = 35
$ws.Range("AB59").Value = 0
$ws.Range("AC59").Value = 3

$ws.Range("P60").Value = $false
$ws.Range("Q60").Value = $false
$ws.Range("R60").Value = $true
$ws.Range("S60").Value = $false
$ws.Range("T60").Value = $false
$ws.Range("U60").Value = $false
$ws.Range("V60").Value = 0
$ws.Range("W60").Value = 0
$ws.Range("X60").Value = 20
$ws.Range("AA60").Value = 20
$ws.Range("AB60").Value = 0
$ws.Range("AC60").Value = 4

# ---------------------------------------------------------------------
# 3. New test sections: TEST 12 (QuickAccel / QuickDecel saturation) and
#    TEST 13 (On pressed with out-of-range speed)
# ---------------------------------------------------------------------

function New-TestHeader($rowNum, $text) {
    $rng = $ws.Range("A" + $rowNum + ":I" + $rowNum)
    $rng.Merge()
    $rng.Font.Bold = $true
    $rng.HorizontalAlignment = -4108
    $ws.Range("A" + $rowNum).Value = $text
}

function New-ColumnHeaders($rowNum) {
    $ws.Range("A$rowNum").Value = "On"
    $ws.Range("B$rowNum").Value = "Off "
    $ws.Range("C$rowNum").Value = "Resume "
    $ws.Range("D$rowNum").Value = "Set"
    $ws.Range("E$rowNum").Value = "QuickAccel"
    $ws.Range("F$rowNum").Value = "QuickDecel"
    $ws.Range("G$rowNum").Value = " Accel"
    $ws.Range("H$rowNum").Value = "Break"
    $ws.Range("I$rowNum").Value = "Speed"
    $ws.Range("L$rowNum").Value = "CruiseSpeed "
    $ws.Range("M$rowNum").Value = "ThrottleCmd "
    $ws.Range("N$rowNum").Value = "CruiseState"
}

function New-DataRow($rowNum, $on, $off, $resume, $set, $quickAccel, $quickDecel, $accel, $brk, $speed, $cruiseSpeed, $throttle, $state) {
    $ws.Range("A$rowNum").Value = $on
    $ws.Range("B$rowNum").Value = $off
    $ws.Range("C$rowNum").Value = $resume
    $ws.Range("D$rowNum").Value = $set
    $ws.Range("E$rowNum").Value = $quickAccel
    $ws.Range("F$rowNum").Value = $quickDecel
    $ws.Range("G$rowNum").Value = $accel
    $ws.Range("H$rowNum").Value = $brk
    $ws.Range("I$rowNum").Value = $speed
    $ws.Range("L$rowNum").Value = $cruiseSpeed
    $ws.Range("M$rowNum").Value = $throttle
    $ws.Range("N$rowNum").Value = $state
}

# --- TEST 12 - QuickAccel button pressed and speed limited to speed max (96-101)
New-TestHeader 96 "TEST 12 - QuickAccel button pressed and speed limited to speed max"
$ws.Range("A97").Value = " "
New-ColumnHeaders 98
New-DataRow 99  $true  $false $false $false $false $false 0 0 35  35  0 2
New-DataRow 100 $false $false $false $false $false $false 0 0 149 149 0 2
New-DataRow 101 $false $false $false $false $true  $false 0 0 149 150 0 2

# --- TEST 12 - QuickDecel button pressed and speed limited to speed min (103-108)
New-TestHeader 103 "TEST 12 - QuickDecel button pressed and speed limited to speed min"
$ws.Range("A104").Value = " "
New-ColumnHeaders 105
New-DataRow 106 $true  $false $false $false $false $false 0 0 35 35 0 2
New-DataRow 107 $false $false $false $false $false $false 0 0 31 31 0 2
New-DataRow 108 $false $false $false $false $false $true  0 0 31 30 0 2

# --- TEST 13 - On pressed with speed above speed max (110-114)
New-TestHeader 110 "TEST 13 - On pressed with speed above speed max"
$ws.Range("A111").Value = " "
New-ColumnHeaders 112
New-DataRow 113 $false $false $false $false $false $false 0 0 160 160 0 1
New-DataRow 114 $true  $false $false $false $false $false 0 0 160 150 0 2

# --- TEST 13 - On pressed with speed below speed min (116-120)
New-TestHeader 116 "TEST 13 - On pressed with speed below speed min"
$ws.Range("A117").Value = " "
New-ColumnHeaders 118
New-DataRow 119 $false $false $false $false $false $false 0 0 20 20 0 1
New-DataRow 120 $true  $false $false $false $false $false 0 0 20 30 0 2

# ---------------------------------------------------------------------
# 4. Update the saved view/selection
# ---------------------------------------------------------------------
$ws.Range("N113").Select()
